$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.119.21'
$ws.Range('E2').Value = '  -0.49%  '

# Row 3
$ws.Range('D3').Value = '2.441.76'
$ws.Range('E3').Value = '  +0.53%  '

# Row 4
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').Value = "'579.80"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.10%  '

# Row 6
$ws.Range('D6').Value = "'143.04"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.83%  '

# Row 7
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').Value = "'0.530"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.34%  '

# Row 9
$ws.Range('D9').Value = '2.438.83'
$ws.Range('E9').Value = '  +0.66%  '

# Row 10
$ws.Range('E10').Value = '  -2.58%  '

# Row 11
$ws.Range('E11').Value = '  +2.66%  '

# Row 12
$ws.Range('D12').Value = "'5.19"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.38%  '

# Row 13
$ws.Range('D13').Value = "'0.343"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.44%  '

# Row 14
$ws.Range('D14').Value = "'26.36"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.08%  '

# Row 15
$ws.Range('D15').Value = "'0.0000171"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.94%  '

# Row 16
$ws.Range('D16').Value = '2.853.19'
$ws.Range('E16').Value = '  -0.24%  '

# Row 17
$ws.Range('D17').Value = '61.986.75'
$ws.Range('E17').Value = '  -0.70%  '

# Row 18
$ws.Range('D18').Value = '2.434.35'
$ws.Range('E18').Value = '  +0.37%  '

# Row 19
$ws.Range('D19').Value = "'10.82"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.10%  '

# Row 20
$ws.Range('D20').Value = "'7.12"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.56%  '

# Row 21
$ws.Range('D21').Value = "'327.73"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.86%  '

# Row 22
$ws.Range('D22').Value = "'4.09"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.63%  '

# Row 23
$ws.Range('E23').Value = '  -3.01%  '

# Row 24
$ws.Range('E24').Value = '  +0.20%  '

# Row 25
$ws.Range('D25').Value = "'65.77"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.53%  '

# Row 26
$ws.Range('D26').Value = "'9.33"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +7.06%  '

# Row 27
$ws.Range('D27').Value = "'613.89"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.11%  '

# Row 28
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.566.63'
$ws.Range('E28').Value = '  -0.52%  '

# Row 29
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0948'
$ws.Range('E29').Value = '  -4.93%  '

# Row 30
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.36%  '

# Row 31
$ws.Range('D31').Value = "'1.43"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -3.67%  '

# Row 32
$ws.Range('D32').Value = "'7.96"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.71%  '

# Row 33
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = "'1.88"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.05%  '

# Row 34
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = "'0.140"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.85%  '

# Row 35
$ws.Range('D35').Value = "'4.88"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.21%  '

# Row 36
$ws.Range('E36').Value = '  +0.42%  '

# Row 37
$ws.Range('E37').Value = '  -4.87%  '

# Row 38
$ws.Range('D38').Value = "'0.375"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.01%  '

# Row 39
$ws.Range('D39').Value = "'149.09"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.14%  '

# Row 40
$ws.Range('D40').Value = "'5.29"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.25%  '

# Row 41
$ws.Range('D41').Value = "'18.29"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.88%  '

# Row 42
$ws.Range('D42').Value = "'1.73"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.06%  '

# Row 43
$ws.Range('D43').Value = "'42.50"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.71%  '

# Row 44
$ws.Range('E44').Value = '  +0.00%  '

# Row 45
$ws.Range('D45').Value = "'2.45"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.42%  '

# Row 46
$ws.Range('D46').Value = "'142.66"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.73%  '

# Row 47
$ws.Range('D47').Value = "'3.62"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.56%  '

# Row 48
$ws.Range('D48').Value = "'0.604"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.90%  '

# Row 49
$ws.Range('D49').Value = "'0.0522"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.89%  '

# Row 50
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = "'19.42"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -6.23%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0238'
$ws.Range('E51').Value = '  +11.41%  '
